# Auto-generated Excel COM-interop script
# Applies 2022-11-14 data update to violent-crime-full-year.xlsx
# For each affected worksheet, updates the 2022 (column I) values
# for the crime-category / neighborhood rows listed in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 6418
$ws.Range('I3').Value = 6698
$ws.Range('I4').Value = 1530
$ws.Range('I5').Value = 625
$ws.Range('I6').Value = 7638
$ws.Range('I7').Value = 22909

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I5').Value = 19
$ws.Range('I6').Value = 212
$ws.Range('I7').Value = 719

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I3').Value = 325
$ws.Range('I5').Value = 26
$ws.Range('I6').Value = 267
$ws.Range('I7').Value = 881

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('I2').Value = 71
$ws.Range('I6').Value = 93
$ws.Range('I7').Value = 235

$ws = $wb.Worksheets.Item('New City')
$ws.Range('I2').Value = 176
$ws.Range('I3').Value = 157
$ws.Range('I7').Value = 528

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I4').Value = 94
$ws.Range('I7').Value = 728
$ws.Range('I8').Value = 1376
$ws.Range('I9').Value = 114
$ws.Range('I10').Value = 161
$ws.Range('I11').Value = 348
$ws.Range('I12').Value = 56
$ws.Range('I15').Value = 266
$ws.Range('I18').Value = 172
$ws.Range('I19').Value = 643
$ws.Range('I20').Value = 568
$ws.Range('I27').Value = 201
$ws.Range('I29').Value = 1390
$ws.Range('I31').Value = 235
$ws.Range('I33').Value = 1032
$ws.Range('I36').Value = 312
$ws.Range('I37').Value = 719
$ws.Range('I40').Value = 41
$ws.Range('I42').Value = 821
$ws.Range('I43').Value = 199
$ws.Range('I44').Value = 172
$ws.Range('I49').Value = 153
$ws.Range('I51').Value = 273
$ws.Range('I52').Value = 493
$ws.Range('I54').Value = 466
$ws.Range('I55').Value = 263
$ws.Range('I60').Value = 127
$ws.Range('I63').Value = 73
$ws.Range('I64').Value = 187
$ws.Range('I65').Value = 528
$ws.Range('I67').Value = 881
$ws.Range('I71').Value = 65
$ws.Range('I72').Value = 91
$ws.Range('I74').Value = 38
$ws.Range('I76').Value = 330
$ws.Range('I79').Value = 651
$ws.Range('I85').Value = 1030
$ws.Range('I90').Value = 295
$ws.Range('I92').Value = 69
$ws.Range('I94').Value = 235
$ws.Range('I95').Value = 347
$ws.Range('I101').Value = 22909

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('I2').Value = 122
$ws.Range('I7').Value = 347

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I2').Value = 234
$ws.Range('I6').Value = 331
$ws.Range('I7').Value = 1032

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('I2').Value = 28
$ws.Range('I6').Value = 92
$ws.Range('I7').Value = 153

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I2').Value = 100
$ws.Range('I6').Value = 224
$ws.Range('I7').Value = 466

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 408
$ws.Range('I3').Value = 479
$ws.Range('I6').Value = 386
$ws.Range('I7').Value = 1390

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I2').Value = 214
$ws.Range('I3').Value = 190
$ws.Range('I7').Value = 643

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('I2').Value = 54
$ws.Range('I6').Value = 50
$ws.Range('I7').Value = 172

$ws = $wb.Worksheets.Item('River North')
$ws.Range('I4').Value = 34
$ws.Range('I7').Value = 330

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I6').Value = 261
$ws.Range('I7').Value = 1030

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I2').Value = 201
$ws.Range('I3').Value = 254
$ws.Range('I6').Value = 284
$ws.Range('I7').Value = 821

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('I3').Value = 33
$ws.Range('I7').Value = 161

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('I2').Value = 79
$ws.Range('I3').Value = 83
$ws.Range('I5').Value = 7
$ws.Range('I7').Value = 263

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('I3').Value = 88
$ws.Range('I6').Value = 65

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I6').Value = 191
$ws.Range('I7').Value = 651

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('I2').Value = 53
$ws.Range('I7').Value = 187

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('I6').Value = 198
$ws.Range('I7').Value = 568

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('I3').Value = 39
$ws.Range('I7').Value = 172

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('I2').Value = 89
$ws.Range('I7').Value = 312

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I2').Value = 133
$ws.Range('I7').Value = 493

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I6').Value = 136
$ws.Range('I7').Value = 235

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('I3').Value = 62
$ws.Range('I7').Value = 266

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('I2').Value = 140
$ws.Range('I3').Value = 72
$ws.Range('I7').Value = 348

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('I2').Value = 37
$ws.Range('I7').Value = 114

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('I6').Value = 29
$ws.Range('I7').Value = 69

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 412
$ws.Range('I7').Value = 1376

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('I4').Value = 26
$ws.Range('I7').Value = 201

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('I2').Value = 96
$ws.Range('I7').Value = 295

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('I6').Value = 110
$ws.Range('I7').Value = 273

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('I3').Value = 33
$ws.Range('I7').Value = 127

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('I3').Value = 34
$ws.Range('I7').Value = 199

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('I4').Value = 4
$ws.Range('I7').Value = 65

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('I6').Value = 43
$ws.Range('I7').Value = 91

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range('I2').Value = 15
$ws.Range('I7').Value = 41

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('I2').Value = 241
$ws.Range('I3').Value = 224
$ws.Range('I7').Value = 728

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('I3').Value = 25
$ws.Range('I7').Value = 94

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('I2').Value = 12
$ws.Range('I7').Value = 56

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range('I3').Value = 2
$ws.Range('I6').Value = 26
$ws.Range('I7').Value = 38
